$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 51, shifting existing rows 51-106 down to 53-108
$ws.Range("A51:A52").EntireRow.Insert()

# New row 51 data
$ws.Cells.Item(51,1).Value = 7
$ws.Cells.Item(51,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(51,3).Value = "Ñuble"
$ws.Cells.Item(51,4).Value = 45128
$ws.Cells.Item(51,5).Value = 16
$ws.Cells.Item(51,6).Value = 100112013
$ws.Cells.Item(51,7).Value = "Alcachofa"
$ws.Cells.Item(51,8).Value = "Argentina(o)"
$ws.Cells.Item(51,9).Value = "Primera"
$ws.Cells.Item(51,10).Value = 50
$ws.Cells.Item(51,11).Value = 17000
$ws.Cells.Item(51,12).Value = 17000
$ws.Cells.Item(51,13).Value = 17000
$ws.Cells.Item(51,14).Value = "$/caja 50 unidades"
$ws.Cells.Item(51,15).Value = "Provincia de Limarí"
$ws.Cells.Item(51,16).Value = 340
$ws.Cells.Item(51,17).Value = 50
$ws.Cells.Item(51,18).Value = "Hortaliza"

# New row 52 data
$ws.Cells.Item(52,1).Value = 7
$ws.Cells.Item(52,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(52,3).Value = "Ñuble"
$ws.Cells.Item(52,4).Value = 45128
$ws.Cells.Item(52,5).Value = 16
$ws.Cells.Item(52,6).Value = 100112013
$ws.Cells.Item(52,7).Value = "Alcachofa"
$ws.Cells.Item(52,8).Value = "Española"
$ws.Cells.Item(52,9).Value = "Primera"
$ws.Cells.Item(52,10).Value = 30
$ws.Cells.Item(52,11).Value = 17000
$ws.Cells.Item(52,12).Value = 17000
$ws.Cells.Item(52,13).Value = 17000
$ws.Cells.Item(52,14).Value = "$/caja 30 unidades"
$ws.Cells.Item(52,15).Value = "Provincia de Limarí"
$ws.Cells.Item(52,16).Value = 567
$ws.Cells.Item(52,17).Value = 30
$ws.Cells.Item(52,18).Value = "Hortaliza"
